$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewSemTests")

# New test-description / path strings, written in the same order the shared
# strings were first introduced (keeps the sharedStrings table identical to
# the authored one): rows 14-19 were drafted first (out of strict row/column
# order), then row 13 was rewritten, with a couple of the B-column
# descriptions for rows 18/19 filled in last.
$ws.Range("B14").Value = 'Integration: P semantics test: one machine, "push" with explicit "pop"'
$ws.Range("A15").Value = 'SEM_OneMachine_13\PushTransInheritance.p'
$ws.Range("B15").Value = 'Integration: P semantics test: one machine, "push" transition, action inherited by the pushed state'
$ws.Range("A14").Value = 'SEM_OneMachine_12\PushExplicitPop.2'
$ws.Range("B16").Value = 'Integration: P semantics test: one machine, "goto" transition, action is not inherited by the destination state'
$ws.Range("A16").Value = 'SEM_OneMachine_14\GotoTransInheritance.p'
$ws.Range("A17").Value = 'SEM_OneMachine_15\ImplicitPopExit.p'
$ws.Range("B17").Value = 'Integration: P semantics test: one machine, exit actions executed upon implicit "pop" '
$ws.Range("A18").Value = 'SEM_OneMachine_16\ExplicitPopExit.p'
$ws.Range("A19").Value = 'SEM_OneMachine_17\PushImplicitPopWithRaise.p'
$ws.Range("A13").Value = 'SEM_OneMachine_11\PushImplicitPopWithSend.p'
$ws.Range("B13").Value = 'Integration: P semantics test: one machine, "push" with implicit "pop" when the unhandled event was sent'
$ws.Range("B18").Value = 'Integration: P semantics test: one machine, exit actions executed upon explicit "pop"'
$ws.Range("B19").Value = 'Integration: P semantics test: one machine, "push" with implicit "pop" when the unhandled event was raised'

# Correct?/Static Error? columns - same "No" / "Yes" pattern as every other
# row in the table.
$ws.Range("C13").Value = "No"
$ws.Range("D13").Value = "Yes"
$ws.Range("C14").Value = "No"
$ws.Range("D14").Value = "Yes"
$ws.Range("C15").Value = "No"
$ws.Range("D15").Value = "Yes"
$ws.Range("C16").Value = "No"
$ws.Range("D16").Value = "Yes"
$ws.Range("C17").Value = "No"
$ws.Range("D17").Value = "Yes"
$ws.Range("C18").Value = "No"
$ws.Range("D18").Value = "Yes"
$ws.Range("C19").Value = "No"
$ws.Range("D19").Value = "Yes"

# "Other features tested" reference for the new Max-Instances related row
$ws.Range("F19").Value = "1.1.1.2. Assert Max Instances of an event"

# Column B grew a bit wider to fit the longer new descriptions
$ws.Columns.Item(2).ColumnWidth = 96.67

# Leave the selection where the author left it after adding the rows
$null = $ws.Range("B22").Select()
